# Trade #40 closed at 2026-02-16 21:29:12 - momentum DOWN +0.000%
# Appends a new "OPEN" trade row (row 12) to the "momentum" sheet,
# mirroring the structure of the existing rows (e.g. row 11).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("momentum")

$row = 12

# --- Plain text / numeric columns -----------------------------------
# These do not look like dates, so a normal .Value assignment is safe
# and keeps the default (unstyled) cell formatting.
$ws.Cells.Item($row, 1).Value = 40              # A: Trade #
$ws.Cells.Item($row, 3).Value = "21:29:12"       # C: Time
$ws.Cells.Item($row, 4).Value = "momentum"       # D: Strategy
$ws.Cells.Item($row, 5).Value = "DOWN"           # E: Side
$ws.Cells.Item($row, 6).Value = 68656.63         # F: Entry Price
$ws.Cells.Item($row, 8).Value = "OPEN"           # H: Status
$ws.Cells.Item($row, 9).Value = 0                # I: P&L %
$ws.Cells.Item($row, 10).Value = 0               # J: P&L $
$ws.Cells.Item($row, 11).Value = 0.9             # K: Confidence
$ws.Cells.Item($row, 12).Value = "Downward momentum: -0.409% over 10 samples"  # L: Entry Reason
$ws.Cells.Item($row, 14).Value = 0               # N: Duration (min)

# --- Date-like text column (B) ---------------------------------------
# "2026-02-16" looks like a date, and a plain .Value/.Formula assignment
# gets auto-converted by Excel into a date serial number (plus a new
# NumberFormat style). To store it as literal text - matching how the
# rest of the sheet stores its "Date" column - build it as a text
# formula and flatten the formula to a static value via copy/paste
# special. That avoids both the date auto-conversion and any new
# styles being introduced.
$bcell = $ws.Cells.Item($row, 2)
$bcell.Formula = '="2026-02-16"'
$bcell.Copy()
$bcell.PasteSpecial(-4163)  # xlPasteValues

# --- Empty columns (G, M) --------------------------------------------
# The source rows for OPEN trades keep an (empty) cell present in the
# Exit Price / Exit Reason columns. A blank .Value assignment removes
# the cell entirely in this engine, so instead "touch" a formatting
# property (a no-op Bold assignment) which keeps the cell present
# without introducing any value/content or new style.
$ws.Cells.Item($row, 7).Font.Bold = $false    # G: Exit Price (empty)
$ws.Cells.Item($row, 13).Font.Bold = $false   # M: Exit Reason (empty)

# --- Keep the sheet's used-range dimension in sync --------------------
$ws.Cells.Item(1, 1).Select()
